$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# These cells hold numeric-looking values that are stored as text (shared
# strings) in the workbook. Setting .Value directly would make Excel
# auto-convert the text to a real number, so we temporarily force a text
# number format, assign the value, then restore the original ("Normal")
# cell style so the cell's formatting ends up unchanged.

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 11: Enterprises density (per 1000 people)
Set-TextValue $ws.Range("C11") "1.65"
Set-TextValue $ws.Range("D11") "56.35"

# Row 12: Employment (% of total)
Set-TextValue $ws.Range("C12") "31.63"
Set-TextValue $ws.Range("D12") "73.23"

# Row 14: Enterprises (% of total)
Set-TextValue $ws.Range("B14") "96.95"
Set-TextValue $ws.Range("C14") "2.92"
Set-TextValue $ws.Range("D14") "99.87"
